# "added the expense overview"
# Expand the Expense sheet from one expense row (Rent) to three rows:
#   row 2: Shopping / 3000 / 2025-12-24
#   row 3: Rent     /  300 / 2025-12-21
#   row 4: transport/ 2000 / 2025-12-10
#
# The date column (C) keeps the existing date-formatted style, so we copy
# that cell's formatting down to the new rows before overwriting values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the date format already applied to C2 ("Rent") down to the two
# new date cells (C3, C4) before any values move around.
$ws.Range("C2").Copy()
$ws.Range("C3:C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: new "Shopping" expense (overwrites the old "Rent" row in place).
$ws.Range("A2").Value = "Shopping"
$ws.Range("B2").Value = 3000
$ws.Range("C2").Value = 46015.22928240741

# Row 3: the original "Rent" expense, moved down with its updated amount
# and date.
$ws.Range("A3").Value = "Rent"
$ws.Range("B3").Value = 300
$ws.Range("C3").Value = 46012.22928240741

# Row 4: new "transport" expense.
$ws.Range("A4").Value = "transport"
$ws.Range("B4").Value = 2000
$ws.Range("C4").Value = 46001.22928240741
